$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D-column price cells are plain text in the workbook (t="inlineStr").
# Assigning a numeric-looking string via .Value would make Excel coerce
# the cell to a number, so force text mode via NumberFormat "@" and restore
# the original cell style afterwards so no stray style index is left behind.
function Set-TextValue($range, $text) {
    $origStyle = $range.Style
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.Style = $origStyle
}

# Row 2
Set-TextValue $ws.Range("D2") "63.212.76"
$ws.Range("E2").Value = "  +0.38%  "

# Row 3
Set-TextValue $ws.Range("D3") "3.083.01"
$ws.Range("E3").Value = "  -1.17%  "

# Row 4
$ws.Range("E4").Value = "  +0.08%  "

# Row 5
Set-TextValue $ws.Range("D5") "580.50"
$ws.Range("E5").Value = "  -0.90%  "

# Row 6
Set-TextValue $ws.Range("D6") "143.81"
$ws.Range("E6").Value = "  -0.59%  "

# Row 7
$ws.Range("E7").Value = "  +0.08%  "

# Row 8
Set-TextValue $ws.Range("D8") "3.078.43"
$ws.Range("E8").Value = "  -1.06%  "

# Row 9
$ws.Range("E9").Value = "  -0.61%  "

# Row 10
Set-TextValue $ws.Range("D10") "0.157"
$ws.Range("E10").Value = "  +4.59%  "

# Row 11
Set-TextValue $ws.Range("D11") "5.60"
$ws.Range("E11").Value = "  -1.83%  "

# Row 12
Set-TextValue $ws.Range("D12") "0.454"
$ws.Range("E12").Value = "  -2.95%  "

# Row 13
$ws.Range("E13").Value = "  -1.37%  "

# Row 14
Set-TextValue $ws.Range("D14") "37.63"
$ws.Range("E14").Value = "  +5.97%  "

# Row 15
$ws.Range("E15").Value = "  -1.23%  "

# Row 16
Set-TextValue $ws.Range("D16") "3.599.79"
$ws.Range("E16").Value = "  -0.85%  "

# Row 17
Set-TextValue $ws.Range("D17") "63.193.08"
$ws.Range("E17").Value = "  +0.50%  "

# Row 18
Set-TextValue $ws.Range("D18") "7.08"
$ws.Range("E18").Value = "  -1.35%  "

# Row 19
Set-TextValue $ws.Range("D19") "3.086.12"
$ws.Range("E19").Value = "  -0.77%  "

# Row 20
Set-TextValue $ws.Range("D20") "459.17"
$ws.Range("E20").Value = "  -1.43%  "

# Row 21
Set-TextValue $ws.Range("D21") "14.12"
$ws.Range("E21").Value = "  +0.30%  "

# Row 22
Set-TextValue $ws.Range("D22") "0.721"
$ws.Range("E22").Value = "  -1.12%  "

# Row 23
Set-TextValue $ws.Range("D23") "7.42"
$ws.Range("E23").Value = "  -1.46%  "

# Row 24
Set-TextValue $ws.Range("D24") "12.92"
$ws.Range("E24").Value = "  -3.40%  "

# Row 25
Set-TextValue $ws.Range("D25") "80.98"
$ws.Range("E25").Value = "  -1.52%  "

# Row 26
Set-TextValue $ws.Range("D26") "2.11"
$ws.Range("E26").Value = "  -3.03%  "

# Row 27
$ws.Range("E27").Value = "  -0.07%  "

# Row 28
Set-TextValue $ws.Range("D28") "8.87"
$ws.Range("E28").Value = "  +6.70%  "

# Row 29
$ws.Range("E29").Value = "  -0.01%  "

# Row 30
$ws.Range("E30").Value = "  -1.01%  "

# Row 31
$ws.Range("E31").Value = "  -2.93%  "

# Row 32
Set-TextValue $ws.Range("D32") "6.79"
$ws.Range("E32").Value = "  -0.83%  "

# Row 33
Set-TextValue $ws.Range("D33") "26.53"
$ws.Range("E33").Value = "  -1.68%  "

# Row 34
$ws.Range("E34").Value = "  -2.03%  "

# Row 35
Set-TextValue $ws.Range("D35") "0.0₃0841"
$ws.Range("E35").Value = "  -0.81%  "

# Row 36
$ws.Range("E36").Value = "  -1.38%  "

# Row 37
Set-TextValue $ws.Range("D37") "2.28"
$ws.Range("E37").Value = "  -3.89%  "

# Row 38
Set-TextValue $ws.Range("D38") "3.32"
$ws.Range("E38").Value = "  +3.17%  "

# Row 39
Set-TextValue $ws.Range("D39") "5.96"
$ws.Range("E39").Value = "  -1.48%  "

# Row 40
Set-TextValue $ws.Range("D40") "50.23"
$ws.Range("E40").Value = "  -1.63%  "

# Row 41
Set-TextValue $ws.Range("D41") "434.16"
$ws.Range("E41").Value = "  +0.42%  "

# Row 42
Set-TextValue $ws.Range("D42") "8.69"
$ws.Range("E42").Value = "  -1.20%  "

# Row 43
Set-TextValue $ws.Range("D43") "0.0367"
$ws.Range("E43").Value = "  -0.52%  "

# Row 44
Set-TextValue $ws.Range("D44") "2.854.51"
$ws.Range("E44").Value = "  -2.74%  "

# Row 45
$ws.Range("E45").Value = "  -2.58%  "

# Row 46
$ws.Range("E46").Value = "  -3.79%  "

# Row 47
Set-TextValue $ws.Range("D47") "35.92"
$ws.Range("E47").Value = "  +1.83%  "

# Row 48
$ws.Range("E48").Value = "  +0.04%  "

# Row 49
Set-TextValue $ws.Range("D49") "123.88"
$ws.Range("E49").Value = "  +0.51%  "

# Row 50
Set-TextValue $ws.Range("D50") "0.109"
$ws.Range("E50").Value = "  -1.78%  "

# Row 51
Set-TextValue $ws.Range("D51") "23.94"
$ws.Range("E51").Value = "  -3.31%  "
